# 盧嘉辰 財產申報表 — sheet #5 (債務) and #6 (事業投資):
# add insurance/claim/debt/investment metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) and
# turn the formerly-duplicated first data row into a proper header row,
# matching the layout already used on sheets 1-4.

$wb = $excel.ActiveWorkbook

# A cell elsewhere in the workbook that already holds the literal text
# "2011-11-18" as a string (sheet 土地, K2). We copy-paste its *value* into
# our new date cells so Excel's smart date parser doesn't turn the literal
# "2011-11-18" we'd otherwise type into a serial date number.
$dateSrc = $wb.Worksheets.Item(1).Range("K2")

# ---------------------------------------------------------------
# Sheet 5: 債務 (Debt)
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Row 1 becomes the header row (previously a duplicate of row 2's data).
$ws5.Range("B1").Value = "species"
$ws5.Range("C1").Value = "debtor"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "register_date"
$ws5.Range("G1").Value = "register_reason"

# New trailing metadata columns H:N on the header row; copy the existing
# bold/bordered header formatting onto them before filling in the labels.
$ws5.Range("G1").Copy() | Out-Null
$ws5.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws5.Range("H1").Value = "property_category"
$ws5.Range("I1").Value = "category"
$ws5.Range("J1").Value = "date"
$ws5.Range("K1").Value = "legislator_name"
$ws5.Range("L1").Value = "legislator_id"
$ws5.Range("M1").Value = "source_file"
$ws5.Range("N1").Value = "index"

# New trailing metadata columns H:N on the data row (row 2).
$ws5.Range("H2").Value = "debt"
$ws5.Range("I2").Value = "normal"
$dateSrc.Copy() | Out-Null
$ws5.Range("J2").PasteSpecial(-4163) | Out-Null
$ws5.Range("K2").Value = "盧嘉辰"
$ws5.Range("L2").Value = 1715
$ws5.Range("M2").Value = "tmp94c1"
$ws5.Range("N2").Value = 96

# ---------------------------------------------------------------
# Sheet 6: 事業投資 (Business investment)
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Row 1 becomes the header row (previously a duplicate of row 2's data).
$ws6.Range("B1").Value = "owner"
$ws6.Range("C1").Value = "company"
$ws6.Range("D1").Value = "address"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"

# New trailing metadata columns H:N on the header row.
$ws6.Range("G1").Copy() | Out-Null
$ws6.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

# New trailing metadata columns H:N on data rows 2 and 3.
$ws6.Range("H2").Value = "investment"
$ws6.Range("I2").Value = "normal"
$dateSrc.Copy() | Out-Null
$ws6.Range("J2").PasteSpecial(-4163) | Out-Null
$ws6.Range("K2").Value = "盧嘉辰"
$ws6.Range("L2").Value = 1715
$ws6.Range("M2").Value = "tmp94c1"
$ws6.Range("N2").Value = 100

$ws6.Range("H3").Value = "investment"
$ws6.Range("I3").Value = "normal"
$dateSrc.Copy() | Out-Null
$ws6.Range("J3").PasteSpecial(-4163) | Out-Null
$ws6.Range("K3").Value = "盧嘉辰"
$ws6.Range("L3").Value = 1715
$ws6.Range("M3").Value = "tmp94c1"
$ws6.Range("N3").Value = 101
